# This script reproduces the commit that adds three new daily "resume scrape"
# worksheets (2023-10-22, 2023-10-21, 2023-10-19) in front of the existing
# workbook, while also splitting what used to be a single "2023-10-04" sheet
# into a renamed/truncated "2023-10-16" sheet plus a brand new "2023-10-04"
# sheet that keeps the original full row of data.

$wb = $excel.ActiveWorkbook

$email = "maheshkulkarni01121@gmail.com"
$phone = "9423627124"

function Set-TextCell($ws, $addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
}

function Set-HeaderRow($ws) {
    $ws.Range("A1").Value = "Email"
    $ws.Range("B1").Value = "Mobile No."
    $ws.Range("C1").Value = "Skills"
}

# ---------------------------------------------------------------------------
# Step 1: the sheet currently named "2023-10-16" (Worksheets.Item(1)) only has
# a single data row (C2 = "R, C, P"). In the target workbook that exact row
# becomes row 2 of the sheet renamed "2023-10-21", with four more rows
# appended below it. Rename + extend it first.
# ---------------------------------------------------------------------------
$sheetOct21 = $wb.Worksheets.Item(1)
$sheetOct21.Name = "2023-10-21"

Set-TextCell $sheetOct21 "C3" "P, C, R"

Set-TextCell $sheetOct21 "A4" $email
Set-TextCell $sheetOct21 "B4" $phone
Set-TextCell $sheetOct21 "C4" "Html, Css, Website, English, Analysis, Training, Engineering, C, Editing, Php, Writing, Database, Video, Networking, Tensorflow, Research, Programming, Technical, Python, Content"

Set-TextCell $sheetOct21 "A5" $email
Set-TextCell $sheetOct21 "B5" $phone
Set-TextCell $sheetOct21 "C5" "Research, Html, Css, Php, Tensorflow, Networking, Database, Content, Website, Analysis, Python, C, Engineering, Training, Technical, English, Writing, Editing, Programming, Video"

Set-TextCell $sheetOct21 "A6" $email
Set-TextCell $sheetOct21 "B6" $phone
Set-TextCell $sheetOct21 "C6" "Writing, C, Css, Research, Editing, Php, Training, Database, Technical, English, Analysis, Tensorflow, Video, Html, Content, Website, Programming, Python, Engineering, Networking"

# ---------------------------------------------------------------------------
# Step 2: brand new sheet "2023-10-22" placed before "2023-10-21", fully
# populated with 6 data rows.
#
# NOTE: Worksheets.Add(Before) re-binds the variable passed in as "Before" to
# the freshly created sheet, so $sheetOct21 must be re-fetched by name
# afterwards to keep referring to the original sheet.
# ---------------------------------------------------------------------------
$sheetOct22 = $wb.Worksheets.Add($sheetOct21)
$sheetOct22.Name = "2023-10-22"
$sheetOct21 = $wb.Worksheets.Item("2023-10-21")
Set-HeaderRow $sheetOct22

Set-TextCell $sheetOct22 "A2" $email
Set-TextCell $sheetOct22 "B2" $phone
Set-TextCell $sheetOct22 "C2" "Website, Editing, Research, Engineering, Python, Programming, Tensorflow, Networking, Training, Html, C, Writing, Database, Analysis, Video, Content, Css, English, Php, Technical"

Set-TextCell $sheetOct22 "A3" $email
Set-TextCell $sheetOct22 "B3" $phone
Set-TextCell $sheetOct22 "C3" "Php, Video, English, Css, Html, Technical, Website, Editing, Analysis, Tensorflow, Networking, Python, Programming, Research, C, Engineering, Training, Writing, Content, Database"

Set-TextCell $sheetOct22 "A4" $email
Set-TextCell $sheetOct22 "B4" $phone
Set-TextCell $sheetOct22 "C4" "Networking, Editing, Website, Programming, C, Html, Css, Database, Technical, Tensorflow, Analysis, English, Php, Research, Python, Engineering, Writing, Video, Training, Content"

Set-TextCell $sheetOct22 "A5" $email
Set-TextCell $sheetOct22 "B5" $phone
Set-TextCell $sheetOct22 "C5" "Html, Database, Writing, Editing, Website, Technical, Content, English, Engineering, Css, C, Video, Training, Python, Php, Analysis, Programming, Tensorflow, Networking, Research"

Set-TextCell $sheetOct22 "A6" $email
Set-TextCell $sheetOct22 "B6" $phone
Set-TextCell $sheetOct22 "C6" "Technical, Php, Tensorflow, English, Programming, Engineering, C, Networking, Training, Website, Html, Python, Editing, Video, Analysis, Database, Research, Css, Content, Writing"

Set-TextCell $sheetOct22 "A7" $email
Set-TextCell $sheetOct22 "B7" $phone
Set-TextCell $sheetOct22 "C7" "Website, Content, Networking, Programming, Css, C, Training, Python, Writing, Research, Html, Engineering, Analysis, Editing, Video, Php, English, Technical, Tensorflow, Database"

# ---------------------------------------------------------------------------
# Step 3: brand new sheet "2023-10-19" placed right after "2023-10-21" (i.e.
# before the old "2023-10-04" sheet), with only the Skills column populated.
# ---------------------------------------------------------------------------
$sheetOct19 = $wb.Worksheets.Add($null, $sheetOct21)
$sheetOct19.Name = "2023-10-19"
Set-HeaderRow $sheetOct19

Set-TextCell $sheetOct19 "C2" "C, P, R"
Set-TextCell $sheetOct19 "C3" "C, R, P"
Set-TextCell $sheetOct19 "C4" "C, R, P"
Set-TextCell $sheetOct19 "C5" "C, R, P"
Set-TextCell $sheetOct19 "C6" "C, R, P"

# ---------------------------------------------------------------------------
# Step 4: the old "2023-10-04" sheet (full mahesh row in A2:C2) needs to be
# duplicated: a brand new "2023-10-04" sheet keeps that original row, while
# the existing sheet object gets renamed to "2023-10-16" and its row 2 is
# trimmed down to just C2 = "R, C, P".
# ---------------------------------------------------------------------------
$sheetOct16 = $wb.Worksheets.Item("2023-10-04")

$sheetOct04 = $wb.Worksheets.Add($null, $sheetOct16)
Set-HeaderRow $sheetOct04
Set-TextCell $sheetOct04 "A2" $email
Set-TextCell $sheetOct04 "B2" $phone
Set-TextCell $sheetOct04 "C2" "Editing, Programming, English, Python, Css, Content, Training, Networking, Research, Database, Engineering, Php, Website, Analysis, Tensorflow, C, Html, Technical, Writing, Video"

$sheetOct16.Name = "2023-10-16"
$sheetOct16.Range("A2").ClearContents()
$sheetOct16.Range("B2").ClearContents()
Set-TextCell $sheetOct16 "C2" "R, C, P"

$sheetOct04.Name = "2023-10-04"

# ---------------------------------------------------------------------------
# The remaining sheets (2023-09-01, 2023-08-31, 2023-08-29, Sheet) are left
# untouched - they only shift position because of the new sheets inserted
# above them.
# ---------------------------------------------------------------------------
